$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Verbs")
Write-Host $ws.Name
Write-Host $wb.Worksheets.Count
